$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply header style (fill+font+alignment) to the new header cells I1:M1 by copying from A1 ---
$ws.Range("A1").Copy()
$ws.Range("I1:M1").PasteSpecial(-4122)

# --- Header row (row 1) text ---
$ws.Range("A1").Value = "Số thứ tự"
$ws.Range("B1").Value = "Mã nhân viên"
$ws.Range("C1").Value = "Số hợp đồng"
$ws.Range("D1").Value = "Số phụ lục"
$ws.Range("E1").Value = "Ngày ký"
$ws.Range("F1").Value = "Ngày hiệu lực"
$ws.Range("G1").Value = "Lương cơ bản"
$ws.Range("H1").Value = "Lương KPI"
$ws.Range("I1").Value = "Phụ cấp ăn trưa"
$ws.Range("J1").Value = "Phụ cấp điện thoại"
$ws.Range("K1").Value = "Phụ cấp khác"
$ws.Range("L1").Value = "Nội dung thay đổi"
$ws.Range("M1").Value = "Ghi chú"

# Update header fill color (navy 4472C4 -> 4F81BD) across the whole header row
$ws.Range("A1:M1").Interior.Color = 12419407
$ws.Range("A1:M1").Interior.PatternColor = 12419407

# Add thin border box around the full header row
$ws.Range("A1:M1").Borders.LineStyle = 1

# --- Data row (row 2) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "NV001"
$ws.Range("C2").Value = "123/2025/HD-MVL"
$ws.Range("D2").Value = "01/2025/PLHD-MVL"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01/01/2025"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "01/01/2025"
$ws.Range("G2").Value = 10000000
$ws.Range("H2").Value = 5000000
$ws.Range("I2").Value = 1000000
$ws.Range("J2").Value = 500000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = "Điều chỉnh lương"
$ws.Range("M2").Value = "'"

# --- Column widths (engine adds +5/6 padding to stored width, so compensate) ---
$ws.Range("D1").ColumnWidth = 20 - (5/6)
$ws.Range("G1").ColumnWidth = 15 - (5/6)
$ws.Range("H1").ColumnWidth = 15 - (5/6)
$ws.Range("I1").ColumnWidth = 15 - (5/6)
$ws.Range("J1").ColumnWidth = 15 - (5/6)
$ws.Range("K1").ColumnWidth = 15 - (5/6)
$ws.Range("L1").ColumnWidth = 40 - (5/6)
$ws.Range("M1").ColumnWidth = 30 - (5/6)
